$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 345 (shift existing data down by 6)
$ws.Range("A345:R350").Insert()

# Populate the newly inserted rows with the new price records
$newData = New-Object "object[,]" 6,18
$newData[0,0] = 10
$newData[0,1] = 'Vega Modelo de Temuco'
$newData[0,2] = 'La Araucanía'
$newData[0,3] = 44551
$newData[0,4] = 9
$newData[0,5] = 100112027
$newData[0,6] = 'Melón'
$newData[0,7] = 'Calameño'
$newData[0,8] = 'Primera'
$newData[0,9] = 950
$newData[0,10] = 1000
$newData[0,11] = 1000
$newData[0,12] = 1000
$newData[0,13] = '$/unidad'
$newData[0,14] = 'Región de O''Higgins'
$newData[0,15] = 1000
$newData[0,16] = 1
$newData[0,17] = 'Hortaliza'
$newData[1,0] = 10
$newData[1,1] = 'Vega Modelo de Temuco'
$newData[1,2] = 'La Araucanía'
$newData[1,3] = 44551
$newData[1,4] = 9
$newData[1,5] = 100112027
$newData[1,6] = 'Melón'
$newData[1,7] = 'Calameño'
$newData[1,8] = 'Segunda'
$newData[1,9] = 1250
$newData[1,10] = 900
$newData[1,11] = 900
$newData[1,12] = 900
$newData[1,13] = '$/unidad'
$newData[1,14] = 'Región de O''Higgins'
$newData[1,15] = 900
$newData[1,16] = 1
$newData[1,17] = 'Hortaliza'
$newData[2,0] = 10
$newData[2,1] = 'Vega Modelo de Temuco'
$newData[2,2] = 'La Araucanía'
$newData[2,3] = 44551
$newData[2,4] = 9
$newData[2,5] = 100112027
$newData[2,6] = 'Melón'
$newData[2,7] = 'Calameño'
$newData[2,8] = 'Tercera'
$newData[2,9] = 900
$newData[2,10] = 800
$newData[2,11] = 800
$newData[2,12] = 800
$newData[2,13] = '$/unidad'
$newData[2,14] = 'Región de O''Higgins'
$newData[2,15] = 800
$newData[2,16] = 1
$newData[2,17] = 'Hortaliza'
$newData[3,0] = 10
$newData[3,1] = 'Vega Modelo de Temuco'
$newData[3,2] = 'La Araucanía'
$newData[3,3] = 44551
$newData[3,4] = 9
$newData[3,5] = 100112027
$newData[3,6] = 'Melón'
$newData[3,7] = 'Tuna'
$newData[3,8] = 'Primera'
$newData[3,9] = 1250
$newData[3,10] = 900
$newData[3,11] = 900
$newData[3,12] = 900
$newData[3,13] = '$/unidad'
$newData[3,14] = 'Región de O''Higgins'
$newData[3,15] = 900
$newData[3,16] = 1
$newData[3,17] = 'Hortaliza'
$newData[4,0] = 10
$newData[4,1] = 'Vega Modelo de Temuco'
$newData[4,2] = 'La Araucanía'
$newData[4,3] = 44551
$newData[4,4] = 9
$newData[4,5] = 100112027
$newData[4,6] = 'Melón'
$newData[4,7] = 'Tuna'
$newData[4,8] = 'Segunda'
$newData[4,9] = 850
$newData[4,10] = 700
$newData[4,11] = 700
$newData[4,12] = 700
$newData[4,13] = '$/unidad'
$newData[4,14] = 'Región de O''Higgins'
$newData[4,15] = 700
$newData[4,16] = 1
$newData[4,17] = 'Hortaliza'
$newData[5,0] = 10
$newData[5,1] = 'Vega Modelo de Temuco'
$newData[5,2] = 'La Araucanía'
$newData[5,3] = 44551
$newData[5,4] = 9
$newData[5,5] = 100112027
$newData[5,6] = 'Melón'
$newData[5,7] = 'Tuna'
$newData[5,8] = 'Tercera'
$newData[5,9] = 750
$newData[5,10] = 500
$newData[5,11] = 500
$newData[5,12] = 500
$newData[5,13] = '$/unidad'
$newData[5,14] = 'Región de O''Higgins'
$newData[5,15] = 500
$newData[5,16] = 1
$newData[5,17] = 'Hortaliza'

$ws.Range("A345:R350").Value = $newData
